$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Update the "last modified" date field text (7/3/2016 1:51:39 PM -> 2:10:26 PM)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("7/3/2016 1:51:39 PM", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "7/3/2016 2:10:26 PM", 2) | Out-Null

function Find-ParagraphIndex($doc, $pattern) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -match $pattern) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 2. Insert a new Caveats bullet "Conversion is a bit slow." just before the
#    "Tables end with \a ..." bullet.
# ---------------------------------------------------------------------------
$idxTables = Find-ParagraphIndex $d "Tables end with \\a"
$pTables = $d.Paragraphs.Item($idxTables)
$rngTables = $pTables.Range
$rngTables.Collapse(1) | Out-Null
$rngTables.InsertParagraphBefore()
$idxTables = Find-ParagraphIndex $d "Tables end with \\a"
$newPara = $d.Paragraphs.Item($idxTables - 1)
$newPara.Range.Text = "Conversion is a bit slow. "

# ---------------------------------------------------------------------------
# 3. Append " Fixme, still doesn't work." to the "header row should be bold"
#    bullet.
# ---------------------------------------------------------------------------
$idxHeader = Find-ParagraphIndex $d "header row should be bold"
$pHeader = $d.Paragraphs.Item($idxHeader)
$rngHeader = $pHeader.Range
$rngHeader.MoveEnd(1, -1) | Out-Null
$endPos = $rngHeader.End
$insertRng = $d.Range($endPos, $endPos)
$insertRng.InsertAfter(" Fixme, still doesn" + [char]0x2019 + "t work.")

# ---------------------------------------------------------------------------
# 4. Add a new final Caveats bullet about the clipboard warning dialog, and
#    relocate the _GoBack bookmark there (it currently sits right before the
#    "Caveats" heading).
# ---------------------------------------------------------------------------
$idxHeader = Find-ParagraphIndex $d "header row should be bold"
$pHeader = $d.Paragraphs.Item($idxHeader)
$rngHeader = $pHeader.Range
$rngHeader.Collapse(0) | Out-Null
$rngHeader.InsertParagraphAfter()
$idxNew = $idxHeader + 1
$newPara2 = $d.Paragraphs.Item($idxNew)
$newPara2.Range.Text = "Big problem that I couldn" + [char]0x2019 + "t solve, dialog when closing word " + [char]0x201C + "You placed a large amount of content on the clipboard" + [char]0x201D + ". Total HEADACHE... reset clipboard to 2 characters, disable alerts, none worked."

# Move the _GoBack bookmark from before "Caveats" to the end of the new bullet.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$newPara2 = $d.Paragraphs.Item($idxNew)
$rngEnd = $newPara2.Range
$rngEnd.MoveEnd(1, -1) | Out-Null
$rngEnd.Collapse(0) | Out-Null
# Paragraphs.Item(...).Range objects don't reliably work as Bookmarks.Add
# anchors this far into the document, so re-wrap as a plain Document.Range.
$bmRng = $d.Range($rngEnd.Start, $rngEnd.End)
$d.Bookmarks.Add("_GoBack", $bmRng) | Out-Null

Write-Output "Done"
